# Add account-information rows for each new role/property (test2..test42),
# mirroring the existing "test1" template row, and append a new "init property"
# module section (the group-header rows that use the TRUE/FALSE style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate column A formatting (text style) from the template row (row 11)
# down through all of the new rows.
$ws.Range("A11").Copy()
$ws.Range("A12:A52").PasteSpecial(-4122)

# Propagate column B "group header" formatting (text style, same as B11) onto
# each new group-header row individually so the literal "1" stays text and not
# a number.
$ws.Range("B11").Copy()
foreach ($r in @(16,21,26,31,33,38,43,45,50)) {
  $ws.Range("B$r").PasteSpecial(-4122)
}

# Fill in the actual values for each new row.
$ws.Range("A12").Value = "test2"
$ws.Range("B12").Value = 1
$ws.Range("A13").Value = "test3"
$ws.Range("B13").Value = 1
$ws.Range("A14").Value = "test4"
$ws.Range("B14").Value = 1
$ws.Range("A15").Value = "test5"
$ws.Range("B15").Value = 1
$ws.Range("A16").Value = "test6"
$ws.Range("B16").Value = "1"
$ws.Range("A17").Value = "test7"
$ws.Range("B17").Value = 1
$ws.Range("A18").Value = "test8"
$ws.Range("B18").Value = 1
$ws.Range("A19").Value = "test9"
$ws.Range("B19").Value = 1
$ws.Range("A20").Value = "test10"
$ws.Range("B20").Value = 1
$ws.Range("A21").Value = "test11"
$ws.Range("B21").Value = "1"
$ws.Range("A22").Value = "test12"
$ws.Range("B22").Value = 1
$ws.Range("A23").Value = "test13"
$ws.Range("B23").Value = 1
$ws.Range("A24").Value = "test14"
$ws.Range("B24").Value = 1
$ws.Range("A25").Value = "test15"
$ws.Range("B25").Value = 1
$ws.Range("A26").Value = "test16"
$ws.Range("B26").Value = "1"
$ws.Range("A27").Value = "test17"
$ws.Range("B27").Value = 1
$ws.Range("A28").Value = "test18"
$ws.Range("B28").Value = 1
$ws.Range("A29").Value = "test19"
$ws.Range("B29").Value = 1
$ws.Range("A30").Value = "test20"
$ws.Range("B30").Value = 1
$ws.Range("A31").Value = "test21"
$ws.Range("B31").Value = "1"
$ws.Range("A32").Value = "test22"
$ws.Range("B32").Value = 1
$ws.Range("A33").Value = "test23"
$ws.Range("B33").Value = "1"
$ws.Range("A34").Value = "test24"
$ws.Range("B34").Value = 1
$ws.Range("A35").Value = "test25"
$ws.Range("B35").Value = 1
$ws.Range("A36").Value = "test26"
$ws.Range("B36").Value = 1
$ws.Range("A37").Value = "test27"
$ws.Range("B37").Value = 1
$ws.Range("A38").Value = "test28"
$ws.Range("B38").Value = "1"
$ws.Range("A39").Value = "test29"
$ws.Range("B39").Value = 1
$ws.Range("A40").Value = "test30"
$ws.Range("B40").Value = 1
$ws.Range("A41").Value = "test31"
$ws.Range("B41").Value = 1
$ws.Range("A42").Value = "test32"
$ws.Range("B42").Value = 1
$ws.Range("A43").Value = "test33"
$ws.Range("B43").Value = "1"
$ws.Range("A44").Value = "test34"
$ws.Range("B44").Value = 1
$ws.Range("A45").Value = "test35"
$ws.Range("B45").Value = "1"
$ws.Range("A46").Value = "test36"
$ws.Range("B46").Value = 1
$ws.Range("A47").Value = "test37"
$ws.Range("B47").Value = 1
$ws.Range("A48").Value = "test38"
$ws.Range("B48").Value = 1
$ws.Range("A49").Value = "test39"
$ws.Range("B49").Value = 1
$ws.Range("A50").Value = "test40"
$ws.Range("B50").Value = "1"
$ws.Range("A51").Value = "test41"
$ws.Range("B51").Value = 1
$ws.Range("A52").Value = "test42"
$ws.Range("B52").Value = 1

# Update the active selection to match the editor's last position.
$ws.Range("B33").Select() | Out-Null
